$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. "(Context ContextKind, LHS ResourceOccurrence, Context Concept
#    Kind, RHS ResourceOccurrence);" paragraph: insert a
#    "/ ResourceOccurrence" after "Context Concept Kind".
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "(Context ContextKind, LHS ResourceOccurrence, Context Concept Kind, RHS ResourceOccurrence);",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(Context ContextKind, LHS ResourceOccurrence, Context Concept Kind / ResourceOccurrence, RHS ResourceOccurrence);",
    2)

# -----------------------------------------------------------------
# 2. "Ontology Matching: state (inferred) equivalence between
#    types, instances, attributes, relationships and values."
#    paragraph: fully replaced with new text.
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Ontology Matching: state (inferred) equivalence between types, instances, attributes, relationships and values.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Apply Context function Kind to input ResourceOccurrences Flux : LHS. Apply Context Concept function Kind to LHS Flux : RHS.",
    2)

# -----------------------------------------------------------------
# 3. "Ontology Merge / Align: state (inferred) equivalence
#    Statements and Upper Ontology mappings between types,
#    instances, attributes, relationships and values." paragraph:
#    its text is replaced, and seven new list paragraphs are
#    inserted right after it (same list-paragraph formatting),
#    ending with the former "Ontology Matching" / "Ontology Merge"
#    sentences (now with " Rules / Grammar." appended).
# -----------------------------------------------------------------
$find = $d.Content.Find
$found = $find.Execute(
    "Ontology Merge / Align: state (inferred) equivalence Statements and Upper Ontology mappings between types, instances, attributes, relationships and values.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $targetRange = $find.Parent.Duplicate
    $para = $targetRange.Paragraphs.First

    # Resolve the 1-based index of this paragraph within the document
    # so we can reliably reach the newly inserted paragraphs afterward
    # (InsertParagraphAfter does not reposition ranges/selections).
    $paraIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Start -eq $para.Range.Start) {
            $paraIndex = $i
            break
        }
    }

    # Replace this paragraph's own text.
    $para.Range.Text = "(Amor, (Pedro, amaA, Maria), amada, Mar" + [char]0x00ED + "a);"

    $newParagraphTexts = @(
        "(Amor, (Pedro, amaA, Maria), Maria, amada);",
        "(Amor, Pedro, (Pedro, amaA, Maria), amante);",
        "(Amor, amante, Pedro, (Pedro, amaA, Maria));",
        "(Amor, Pedro, amante, (Pedro, amaA, Maria));",
        "(etc.: CSPO, Kinds, Statements LHS, Concepts, RHS).",
        "Ontology Matching: state (inferred) equivalence between types, instances, attributes, relationships and values. Rules / Grammar.",
        "Ontology Merge / Align: state (inferred) equivalence Statements and Upper Ontology mappings between types, instances, attributes, relationships and values. Rules / Grammar."
    )

    $idx = $paraIndex
    foreach ($t in $newParagraphTexts) {
        $endRng = $d.Paragraphs($idx).Range.Duplicate
        $endRng.Collapse(0)
        $endRng.InsertParagraphAfter()
        $idx = $idx + 1
        $d.Paragraphs($idx).Range.Text = $t
    }
}
